# Generate Report for Archive
#
# 1. Update the localization "Status" value from "Ready for handoff" to
#    "In Translation" everywhere it appears:
#      - Overview sheet: E2 (zh-cn status) and F2 (de-de status)
#      - zh-cn sheet:     C2 (Status column)
#      - de-de sheet:     C2 (Status column)
#
# 2. Narrow the "Status" column (and the mirrored zh-cn/de-de columns on the
#    Overview sheet) now that the status text is shorter.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
